# Ver 6 - Special Letters and Trustee Sigs
#
# This edit:
#   1) Removes the stray "_GoBack" bookmark that sat at the top of the
#      "As more people insist ..." paragraph.
#   2) Removes the "SET SIGNATURE "DAVE"" field (the fldChar begin /
#      instrText / fldChar end run trio) that followed "Sincerely,".
#   3) Re-creates the "_GoBack" bookmark immediately after "Sincerely,"
#      (i.e. where the user's cursor was when the field was deleted),
#      wrapping it around the existing "BODY" bookmark's closing mark.

$d = $word.ActiveDocument

# --- Step 1: drop the old "_GoBack" bookmark -------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 2: remove the "SET SIGNATURE "DAVE"" field ------------------------
for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $fld = $d.Fields.Item($i)
    if ($fld.Code.Text -like "*SET SIGNATURE*") {
        $fld.Delete()
    }
}

# --- Step 3: rebuild the final ("Sincerely,") paragraph so the new
#     "_GoBack" bookmark lands right after the text and wraps the
#     existing bookmarkEnd (id 0) that closes the "BODY" bookmark. -------
$lastPara = $d.Paragraphs.Last
$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00261C27" w:rsidRPr="00261C27" w:rsidRDefault="00261C27" w:rsidP="001A59D3">
  <w:pPr>
    <w:spacing w:after="200"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:noProof/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
  </w:pPr>
  <w:r w:rsidRPr="00261C27">
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:noProof/>
      <w:sz w:val="21"/>
      <w:szCs w:val="21"/>
    </w:rPr>
    <w:t>Sincerely,</w:t>
  </w:r>
  <w:bookmarkStart w:id="1" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
  <w:bookmarkEnd w:id="1"/>
</w:p>
'@

[void]$lastPara.Range.InsertXML($newParaXml)
